$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 23:03"

# --- Swap adjacent country rows (labels only; numeric data stays put) ---
# Brunei (row 147) / Zambia (row 148)
$ws.Range("A147").Value = "Zambia"
$ws.Range("A148").Value = "Brunei"

# Burundi (row 199) / San Cristobal y Nieves (row 200)
$ws.Range("A199").Value = "San Cristobal y Nieves"
$ws.Range("A200").Value = "Burundi"

# Seychelles (row 205) / Montserrat (row 206)
$ws.Range("A205").Value = "Montserrat"
$ws.Range("A206").Value = "Seychelles"

# --- Update numeric data ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1232604
$ws.Range("C4").Value = 19769
$ws.Range("E4").Value = 961831
$ws.Range("G4").Value = 1917
$ws.Range("H4").Value = 71838

# Row 15: Canada
$ws.Range("B15").Value = 61981
$ws.Range("C15").Value = 1209
$ws.Range("D15").Value = 26661
$ws.Range("E15").Value = 31283
$ws.Range("G15").Value = 183
$ws.Range("H15").Value = 4037

# Row 134: Ruanda
$ws.Range("D134").Value = 129
$ws.Range("E134").Value = 132

# Row 147
$ws.Range("C147").Value = 1
$ws.Range("D147").Value = 92
$ws.Range("E147").Value = 43
$ws.Range("F147").Value = 1
$ws.Range("H147").Value = 3

# Row 148
$ws.Range("B148").Value = 138
$ws.Range("D148").Value = 131
$ws.Range("E148").Value = 6
$ws.Range("F148").Value = 2
$ws.Range("H148").Value = 1

# Row 164: Barbados
$ws.Range("D164").Value = 47
$ws.Range("E164").Value = 28

# Row 199
$ws.Range("D199").Value = 8
$ws.Range("H199").Value = 0

# Row 200
$ws.Range("D200").Value = 7
$ws.Range("H200").Value = 1

# Row 205
$ws.Range("D205").Value = 7
$ws.Range("F205").Value = 1
$ws.Range("H205").Value = 1

# Row 206
$ws.Range("D206").Value = 8
$ws.Range("F206").Value = 0
$ws.Range("H206").Value = 0
